$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New upcoming-match rows appended at the bottom of the table (rows 116-119).
# Formatting (styles for column A / D) is copied from the last existing data
# row so the new rows match the sheet's established look (bold/border id
# style, and the custom date/time number format), then values are filled in.

$newRows = @(
    @{ Row = 116; Id = 114; MatchId = "7802881"; Date = 45444.70833333334; Home = "Forge FC";        Away = "York United FC"; L = 1.8;   M = 3.25; N = 4;   O = 1.85;  P = 3.2; Q = 3.8;  R = -0.5;  S = 1.9;   T = 1.9;   U = 2.5; V = 1.875; W = 1.925 },
    @{ Row = 117; Id = 115; MatchId = "7802946"; Date = 45444.83333333334; Home = "Pacific FC CA";    Away = "Cavalry FC";      L = 2.5;   M = 3.2;  N = 2.5; O = 2.7;   P = 3.2; Q = 2.35; R = 0.25;  S = 1.725; T = 2.075; U = 2.5; V = 2.025; W = 1.775 },
    @{ Row = 118; Id = 116; MatchId = "7802947"; Date = 45445.625;         Home = "Atletico Ottawa";  Away = "HFX Wanderers";   L = 1.615; M = 3.4;  N = 5;   O = 1.666; P = 3.3; Q = 4.5;  R = -0.75; S = 1.9;   T = 1.9;   U = 2.5; V = 2.025; W = 1.775 },
    @{ Row = 119; Id = 117; MatchId = "7803370"; Date = 45445.75;          Home = "Valour FC";        Away = "Vancouver FC";    L = 2.6;   M = 3.2;  N = 2.4; O = 2.6;   P = 3.1; Q = 2.45; R = 0;     S = 1.975; T = 1.825; U = 2.5; V = 1.9;   W = 1.9 }
)

foreach ($nr in $newRows) {
    $r = $nr.Row

    $ws.Range("A115:AD115").Copy() | Out-Null
    $ws.Range("A" + $r + ":AD" + $r).PasteSpecial(-4122) | Out-Null

    # These are fixtures yet to be played: no full/half-time score and no
    # result letter, and the AA:AD post-line-movement columns aren't used.
    $ws.Range("G" + $r + ":K" + $r).ClearContents() | Out-Null
    $ws.Range("AA" + $r + ":AD" + $r).ClearContents() | Out-Null

    $ws.Cells.Item($r, 1).Value = $nr.Id
    $ws.Cells.Item($r, 3).Value = "Canada Premier League"
    $ws.Cells.Item($r, 4).Value = $nr.Date
    $ws.Cells.Item($r, 5).Value = $nr.Home
    $ws.Cells.Item($r, 6).Value = $nr.Away

    # Column B holds the external match id. For these new fixtures the
    # source feed supplies it as text rather than a plain integer, so force
    # text storage (otherwise Excel auto-coerces the numeric-looking string
    # back into a number) and then drop the temporary formatting footprint.
    $bcell = $ws.Cells.Item($r, 2)
    $bcell.NumberFormat = "@"
    $bcell.Value = $nr.MatchId
    $bcell.Style = "Normal"

    $ws.Cells.Item($r, 12).Value = $nr.L
    $ws.Cells.Item($r, 13).Value = $nr.M
    $ws.Cells.Item($r, 14).Value = $nr.N
    $ws.Cells.Item($r, 15).Value = $nr.O
    $ws.Cells.Item($r, 16).Value = $nr.P
    $ws.Cells.Item($r, 17).Value = $nr.Q
    $ws.Cells.Item($r, 18).Value = $nr.R
    $ws.Cells.Item($r, 19).Value = $nr.S
    $ws.Cells.Item($r, 20).Value = $nr.T
    $ws.Cells.Item($r, 21).Value = $nr.U
    $ws.Cells.Item($r, 22).Value = $nr.V
    $ws.Cells.Item($r, 23).Value = $nr.W
    $ws.Cells.Item($r, 24).Value = 0
    $ws.Cells.Item($r, 25).Value = 0
    $ws.Cells.Item($r, 26).Value = 0
}
